$wb = $excel.ActiveWorkbook

# Column F ("想去人数") updates on sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value  = 172
$ws1.Range("F7").Value  = 2509
$ws1.Range("F8").Value  = 158
$ws1.Range("F11").Value = 1512
$ws1.Range("F12").Value = 522
$ws1.Range("F14").Value = 325
$ws1.Range("F18").Value = 205
$ws1.Range("F22").Value = 156
$ws1.Range("F24").Value = 1596
$ws1.Range("F27").Value = 570

# Column F ("想去人数") updates on sheet "全部类型" (same events, offset rows)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value  = 172
$ws4.Range("F8").Value  = 2509
$ws4.Range("F9").Value  = 158
$ws4.Range("F12").Value = 1512
$ws4.Range("F13").Value = 522
$ws4.Range("F15").Value = 325
$ws4.Range("F19").Value = 205
$ws4.Range("F23").Value = 156
$ws4.Range("F25").Value = 1596
$ws4.Range("F28").Value = 570
